{"js": "// Update the GitHub hyperlink in the document: the link text and the\n// hyperlink target both change from the \"Assignment 4 String Matching\"\n// URL to the \"Assignment 2 External Merge Sort\" URL (ExternalMergeSort.java).\n\nconst oldUrl =\n  \"https://github.com/Rohan-Sharma03/AdvancedDSA/tree/master/Assignment%204%20String%20Matching%20Algo/Navie%2C%20Rabin%20Karp%2C%20KMP\";\nconst newUrl =\n  \"https://github.com/Rohan-Sharma03/AdvancedDSA/blob/master/Assignment%202%20External%20Merge%20Sort/ExternalMergeSort.java\";\n\nconst body = context.document.body;\n\n// Locate the run of text inside the existing hyperlink field.\nconst results = body.search(oldUrl, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the existing hyperlink text to replace.\");\n}\n\nconst target = results.items[0];\n\n// Replace the visible text in place (keeps it inside the hyperlink field).\ntarget.insertText(newUrl, Word.InsertLocation.replace);\nawait context.sync();\n\n// Re-locate the (now updated) text and repoint the hyperlink itself at the\n// new address so the field target matches the displayed text.\nconst results2 = body.search(newUrl, { matchCase: true });\nresults2.load(\"items\");\nawait context.sync();\n\nconst updated = results2.items[0];\nupdated.hyperlink = newUrl;\nawait context.sync();\n", "ps1": "# Update the GitHub hyperlink in the document: both the displayed text and\n# the hyperlink target change from the \"Assignment 4 String Matching\" URL\n# to the \"Assignment 2 External Merge Sort\" URL (ExternalMergeSort.java).\n\n$d = $word.ActiveDocument\n\n$newUrl = \"https://github.com/Rohan-Sharma03/AdvancedDSA/blob/master/Assignment%202%20External%20Merge%20Sort/ExternalMergeSort.java\"\n\n$h = $d.Hyperlinks.Item(1)\n$h.TextToDisplay = $newUrl\n$h.Address = $newUrl\n"}
